{"js": "// Update the date heading and the 20x5 practice-problem table with the\n// \"answer key\" values from the next day's worksheet.\n//\n// Strategy: the body starts with a single centered paragraph holding the\n// date, followed by one table (20 rows x 5 columns) of \"a op b = c\" cells.\n// We replace the date paragraph's text in place (preserving its run\n// formatting) and then overwrite the whole table via `table.values =`,\n// which Word's Office.js API applies cell-by-cell (reusing each cell's\n// existing paragraph/run so fonts/sizes stay intact).\n\nconst newDateText = \"2025-03-20 Thursday\";\n\nconst newTableValues = [\n  [\"83-14=69\", \"58+17=75\", \"63-5=58\", \"7+86=93\", \"8+13=21\"],\n  [\"84-39=45\", \"40-12=28\", \"63+8=71\", \"36+35=71\", \"67-38=29\"],\n  [\"93-88=5\", \"18+27=45\", \"14+39=53\", \"20-16=4\", \"76-9=67\"],\n  [\"80-49=31\", \"14+29=43\", \"61-9=52\", \"67+19=86\", \"76-47=29\"],\n  [\"48+19=67\", \"78+4=82\", \"7+84=91\", \"52-46=6\", \"45+26=71\"],\n  [\"83-78=5\", \"75-66=9\", \"73-49=24\", \"92-44=48\", \"54+18=72\"],\n  [\"75-48=27\", \"18+44=62\", \"76-59=17\", \"38+26=64\", \"58-39=19\"],\n  [\"86-38=48\", \"52-19=33\", \"6+39=45\", \"62-45=17\", \"26+39=65\"],\n  [\"50-24=26\", \"3+59=62\", \"18+49=67\", \"7+89=96\", \"38+15=53\"],\n  [\"3+59=62\", \"57+37=94\", \"59+4=63\", \"80-17=63\", \"77-69=8\"],\n  [\"45+9=54\", \"83-75=8\", \"84-28=56\", \"9+6=15\", \"24+7=31\"],\n  [\"37+39=76\", \"37+29=66\", \"19+66=85\", \"78+4=82\", \"27+16=43\"],\n  [\"19+8=27\", \"93-57=36\", \"19+22=41\", \"16+16=32\", \"30-28=2\"],\n  [\"72-4=68\", \"93-24=69\", \"7+69=76\", \"50-42=8\", \"50-15=35\"],\n  [\"91-89=2\", \"19+69=88\", \"66+15=81\", \"23-17=6\", \"28+56=84\"],\n  [\"56-47=9\", \"94-28=66\", \"86-38=48\", \"62-15=47\", \"74-29=45\"],\n  [\"73-64=9\", \"37+35=72\", \"56-28=28\", \"45-27=18\", \"70-15=55\"],\n  [\"92-13=79\", \"31-17=14\", \"52-17=35\", \"38+24=62\", \"9+19=28\"],\n  [\"50-43=7\", \"56+5=61\", \"91-3=88\", \"42+39=81\", \"64-35=29\"],\n  [\"91-18=73\", \"54+8=62\", \"53-6=47\", \"72-3=69\", \"45-36=9\"],\n];\n\n// 1. Update the date heading (first paragraph of the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n\nif (dateParagraph.text.trim() !== newDateText) {\n  dateParagraph.insertText(newDateText, Word.InsertLocation.replace);\n}\n\n// 2. Update the table of equations, row by row / cell by cell.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newTableValues;\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 20x5 practice-problem table with the\n# \"answer key\" values from the next day's worksheet.\n#\n# The document body is: one centered paragraph holding the date, followed\n# by a single table (20 rows x 5 columns) of \"a op b = c\" cells. We update\n# the date paragraph's Range.Text in place (keeps its run formatting) and\n# then walk the table Cell-by-Cell (row-major, matching reading order) so\n# each of the two repeated \"old\" equations (\"53+39=92\" and \"25+19=44\")\n# gets the correct positional replacement instead of a global text swap.\n\n$d = $word.ActiveDocument\n\n# 1. Update the date heading (first paragraph of the body).\n$d.Paragraphs.Item(1).Range.Text = \"2025-03-20 Thursday\"\n\n# 2. Update the table of equations, row by row / cell by cell.\n$newValues = @(\n    @(\"83-14=69\", \"58+17=75\", \"63-5=58\", \"7+86=93\", \"8+13=21\"),\n    @(\"84-39=45\", \"40-12=28\", \"63+8=71\", \"36+35=71\", \"67-38=29\"),\n    @(\"93-88=5\", \"18+27=45\", \"14+39=53\", \"20-16=4\", \"76-9=67\"),\n    @(\"80-49=31\", \"14+29=43\", \"61-9=52\", \"67+19=86\", \"76-47=29\"),\n    @(\"48+19=67\", \"78+4=82\", \"7+84=91\", \"52-46=6\", \"45+26=71\"),\n    @(\"83-78=5\", \"75-66=9\", \"73-49=24\", \"92-44=48\", \"54+18=72\"),\n    @(\"75-48=27\", \"18+44=62\", \"76-59=17\", \"38+26=64\", \"58-39=19\"),\n    @(\"86-38=48\", \"52-19=33\", \"6+39=45\", \"62-45=17\", \"26+39=65\"),\n    @(\"50-24=26\", \"3+59=62\", \"18+49=67\", \"7+89=96\", \"38+15=53\"),\n    @(\"3+59=62\", \"57+37=94\", \"59+4=63\", \"80-17=63\", \"77-69=8\"),\n    @(\"45+9=54\", \"83-75=8\", \"84-28=56\", \"9+6=15\", \"24+7=31\"),\n    @(\"37+39=76\", \"37+29=66\", \"19+66=85\", \"78+4=82\", \"27+16=43\"),\n    @(\"19+8=27\", \"93-57=36\", \"19+22=41\", \"16+16=32\", \"30-28=2\"),\n    @(\"72-4=68\", \"93-24=69\", \"7+69=76\", \"50-42=8\", \"50-15=35\"),\n    @(\"91-89=2\", \"19+69=88\", \"66+15=81\", \"23-17=6\", \"28+56=84\"),\n    @(\"56-47=9\", \"94-28=66\", \"86-38=48\", \"62-15=47\", \"74-29=45\"),\n    @(\"73-64=9\", \"37+35=72\", \"56-28=28\", \"45-27=18\", \"70-15=55\"),\n    @(\"92-13=79\", \"31-17=14\", \"52-17=35\", \"38+24=62\", \"9+19=28\"),\n    @(\"50-43=7\", \"56+5=61\", \"91-3=88\", \"42+39=81\", \"64-35=29\"),\n    @(\"91-18=73\", \"54+8=62\", \"53-6=47\", \"72-3=69\", \"45-36=9\")\n)\n\n$tbl = $d.Tables.Item(1)\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $rowValues.Count; $c++) {\n        $tbl.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
